$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at the very top. This pushes the existing header
# row (row 1 -> row 2) and every data row below it down by one, carrying
# their formatting (including the bold/centered/bordered style on the
# header row) along for the ride.
$ws.Rows.Item(1).Insert()

# Give the new row 1 the same formatting as the header row that is now on
# row 2 (bold, centered, bordered "s=1" style) before filling it with its
# new content.
$ws.Range("A2:P2").Copy()
$ws.Range("A1:P1").PasteSpecial(-4122)

# New row 1 becomes a simple numeric column index: 0, 1, 2, ... 15.
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
$ws.Range("K1").Value = 10
$ws.Range("L1").Value = 11
$ws.Range("M1").Value = 12
$ws.Range("N1").Value = 13
$ws.Range("O1").Value = 14
$ws.Range("P1").Value = 15

# Row 2 (the former header row, now shifted down) keeps its text labels,
# but the "Pkg." sub-column (M, already blank) and the two API-only
# columns (thread_size / material_surface) are cleared out on this row.
$ws.Range("M2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("P2").ClearContents()
